# Updated symbol list on Sun Jan 29 06:38:19 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.76"
$ws.Range("E2").Value = "'-0.60%"

$ws.Range("D3").Value = "'39.77"
$ws.Range("E3").Value = "'2.79%"

$ws.Range("D4").Value = "'5.139"
$ws.Range("E4").Value = "'0.25%"

$ws.Range("D5").Value = "'0.08134"
$ws.Range("E5").Value = "'-0.56%"

$ws.Range("D6").Value = "'1.947"
$ws.Range("E6").Value = "'-3.11%"

$ws.Range("D7").Value = "'8.155"
$ws.Range("E7").Value = "'2.90%"

$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9273"
$ws.Range("E8").Value = "'-0.51%"

$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1426"
$ws.Range("E9").Value = "'1.42%"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1933"
$ws.Range("E10").Value = "'-0.82%"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09097"
$ws.Range("E11").Value = "'-0.18%"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03531"
$ws.Range("E12").Value = "'1.78%"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09813"
$ws.Range("E13").Value = "'-0.41%"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001395"
$ws.Range("E14").Value = "'-1.45%"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005833"
$ws.Range("E15").Value = "'-1.38%"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.924"
$ws.Range("E16").Value = "'9.91%"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.229"
$ws.Range("E17").Value = "'0.76%"

$ws.Range("E18").Value = "'-1.62%"

$ws.Range("E19").Value = "'-0.14%"

$ws.Range("D20").Value = "'0.1313"
$ws.Range("E20").Value = "'-0.15%"

$ws.Range("D21").Value = "'4.720"
$ws.Range("E21").Value = "'-1.82%"

$ws.Range("E22").Value = "'-0.67%"

$ws.Range("D23").Value = "'0.04376"
$ws.Range("E23").Value = "'-2.01%"

$ws.Range("D24").Value = "'0.001230"
$ws.Range("E24").Value = "'-0.79%"

$ws.Range("D25").Value = "'0.004382"
$ws.Range("E25").Value = "'5.00%"

$ws.Range("E26").Value = "'-0.17%"

$ws.Range("D27").Value = "'0.0004004"
$ws.Range("E27").Value = "'-9.99%"

$ws.Range("E39").Value = "'-3.39%"

$ws.Range("D40").Value = "'0.05098"
$ws.Range("E40").Value = "'-1.67%"

$ws.Range("D41").Value = "'0.007416"
$ws.Range("E41").Value = "'-0.79%"

$ws.Range("D42").Value = "'0.009897"
$ws.Range("E42").Value = "'-1.01%"

$ws.Range("E44").Value = "'-0.17%"

$ws.Range("D45").Value = "'0.009558"
$ws.Range("E45").Value = "'-2.13%"

$ws.Range("D46").Value = "'0.00006371"
$ws.Range("E46").Value = "'0.48%"

$ws.Range("E47").Value = "'-0.16%"

$ws.Range("E49").Value = "'-18.85%"

$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.16%"

$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.16%"
